$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Ngũ hành bản Mệnh sinh Ngũ hành Thiên Mã"
$ws.Range("B8").Value = "Bạn phải lao tâm, vất vả lo nghĩ về những thay đổi, bạn rất năng động và hay phải đi xa rất nhiều."

$ws.Range("A9").Value = "Ngũ hành bản Mệnh khắc Ngũ hành Thiên Mã"
$ws.Range("B9").Value = "Bạn "

$ws.Range("A10").Value = "Ngũ hành Thiên Mã sinh Ngũ hành bản Mệnh"

$ws.Range("A11").Value = "Ngũ hành Thiên Mã khắc Ngũ hành bản Mệnh"

$ws.Range("A12").Value = "Ngũ hành Thiên Mã đồng hành cùng Ngũ hành bản Mệnh"
$ws.Range("B12").Value = "Bạn đi lai, di chuyển an toàn, khi bạn tích cực và năng động tính toán lo lắng công việc thì mọi chuyện hanh thông. "

$ws.Range("B12").Select()
